# The deck originally held 4 slides (SlideID 256, 257, 258, 259). The
# author pruned it down to just the one slide with SlideID 258 (the
# "timeline" slide), dropping the other three (two real slides plus a
# trailing blank one) entirely.
$p = $ppt.ActivePresentation

for ($i = $p.Slides.Count; $i -ge 1; $i--) {
    $sl = $p.Slides.Item($i)
    if ($sl.SlideID -ne 258) {
        $sl.Delete()
    }
}

# On the surviving slide, every remaining top-level shape (the two
# background pictures, the nested "Group 11", the legend picture, the
# four leader-line connectors, and the title picture) got wrapped into one
# new outer group, then that group was nudged to a new position.
$s = $p.Slides.Item(1)

$count = $s.Shapes.Count
$ids = @()
for ($i = 1; $i -le $count; $i++) { $ids += $i }

$grp = $s.Shapes.Range($ids).Group()
$grp.Name = "Group 20"

# Move the freshly-created group to match the authored layout
# (2341566, 983993) EMU, expressed in points (1 pt = 12700 EMU) with
# enough precision to survive the Shape.Left/.Top single-precision
# round-trip exactly.
$grp.Left = 184.3753
$grp.Top = 77.4798
